# Remove the "BCR" (Maximum construction density) and "FAR" (Coefficient of
# land use) rows from the land-use attribute table — they were the last two
# rows of the second table in the document.
$d = $word.ActiveDocument

$t = $d.Tables.Item(2)

# Delete from the bottom up so row indices of the rows still to be removed
# don't shift while we work.
$t.Rows.Item(7).Delete()
$t.Rows.Item(6).Delete()
